# Generate Report for Archive
#
# 1. The localization status for the (only) two source files has moved on
#    from "Ready for handoff" to "In Translation". That status string is
#    shown on the Overview sheet (columns E/F, one per locale) and on each
#    locale's own sheet (column C, "Status").
# 2. The "Status" columns got narrower now that the new status text is
#    shorter than the old one (autosized by the report generator).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Update every cell holding the old status text -----------------
# xlWhole (1) so we only touch cells that equal the old status exactly,
# not e.g. "Latest Handoff File" which merely contains a similar word.
$xlWhole = 1
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    [void]$used.Replace($oldStatus, $newStatus, $xlWhole)
}

# --- 2. Narrow the Status columns to match the shorter text ------------
# Overview!E:F ("zh-cn" / "de-de" status columns)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").ColumnWidth = 12.5

# Each locale sheet's Status column (column C)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
